# Generate Report for Handback
# This script updates the localization-status report to reflect that the
# handback files are now in sync with en-US (i.e. the handback has been
# generated/accepted), clearing the previous "version mismatch" error and
# refreshing the handback timestamps.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status columns (zh-cn / de-de) ---
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$ws2.Range("C2").Value = $newStatus                 # Status
$ws2.Range("K2").Value = "2016-09-04 02:53:46"      # Latest Handback DateTime
$ws2.Range("P2").Value = ""                         # Error Detail (no more mismatch)

# --- de-de sheet ---
$ws3.Range("C2").Value = $newStatus                 # Status
$ws3.Range("K2").Value = "2016-09-04 02:53:53"      # Latest Handback DateTime
$ws3.Range("P2").Value = ""                         # Error Detail (no more mismatch)

# --- Column width adjustments (widen status-related columns, narrow now-empty
#     Error Detail column) ---
$ws1.Columns.Item(5).ColumnWidth = 29.144371396019366   # E: zh-cn status column
$ws1.Columns.Item(6).ColumnWidth = 29.144371396019366   # F: de-de status column

$ws2.Columns.Item(3).ColumnWidth = 29.144371396019366   # C: Status
$ws2.Columns.Item(16).ColumnWidth = 12.913719540550566  # P: Error Detail

$ws3.Columns.Item(3).ColumnWidth = 29.144371396019366   # C: Status
$ws3.Columns.Item(16).ColumnWidth = 12.913719540550566  # P: Error Detail
